$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 28511.555
$ws.Range("I10").Value = 204
$ws.Range("J10").Value = 32050
$ws.Range("K10").Value = 204
$ws.Range("L10").Value = 32050
$ws.Range("M10").Value = 89
$ws.Range("N10").Value = -32636
# Row 13
$ws.Range("H13").Value = 995
$ws.Range("J13").Value = 995
$ws.Range("L13").Value = 995
$ws.Range("N13").Value = -1333
# Row 141
$ws.Range("H141").Value = 2866
$ws.Range("I141").Value = 2499.5
$ws.Range("K141").Value = 7498.5
$ws.Range("M141").Value = -2318.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 19733.334
$ws.Range("J8").Value = 19733.334
$ws.Range("L8").Value = 19733.334
$ws.Range("N8").Value = -20021.334
# Row 10
$ws.Range("H10").Value = 1200
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1200
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -1540
# Row 11
$ws.Range("H11").Value = 23480
$ws.Range("J11").Value = 23480
$ws.Range("L11").Value = 23480
$ws.Range("N11").Value = -23768
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
# Row 41
$ws.Range("H41").Value = 3777.4
$ws.Range("I41").Value = 2499.75
$ws.Range("J41").Value = 8888
$ws.Range("K41").Value = 2499.75
$ws.Range("L41").Value = 8888
$ws.Range("M41").Value = -2085.75
$ws.Range("N41").Value = -9716

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 5686.222
$ws.Range("I5").Value = 46
$ws.Range("J5").Value = 16966.666
$ws.Range("K5").Value = 46
$ws.Range("L5").Value = 16966.666
$ws.Range("M5").Value = 67
$ws.Range("N5").Value = -17192.666
# Row 7
$ws.Range("H7").Value = 3800700
$ws.Range("I7").Value = 4750750
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 4750750
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = -4750637
$ws.Range("N7").Value = -726
# Row 76
$ws.Range("H76").Value = 11500
$ws.Range("J76").Value = 11500
$ws.Range("L76").Value = 11500
$ws.Range("N76").Value = -12130
# Row 79
$ws.Range("H79").Value = 11500
$ws.Range("J79").Value = 11500
$ws.Range("L79").Value = 11500
$ws.Range("N79").Value = -13684
# Row 86
$ws.Range("H86").Value = 4797.294
$ws.Range("I86").Value = 5129.5
$ws.Range("K86").Value = 5129.5
$ws.Range("M86").Value = -4006.5
# Row 89
$ws.Range("H89").Value = 4797.294
$ws.Range("I89").Value = 5129.5
$ws.Range("K89").Value = 25647.5
$ws.Range("M89").Value = -20031.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 434
$ws.Range("I2").Value = 502
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 502
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -389
$ws.Range("N2").Value = -626

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 562.6667
$ws.Range("I97").Value = 601.5
$ws.Range("J97").Value = 485
$ws.Range("K97").Value = 1804.5
$ws.Range("L97").Value = 1455
$ws.Range("M97").Value = -1308.5
$ws.Range("N97").Value = -2447

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 5578.8423
$ws.Range("J3").Value = 3666
$ws.Range("L3").Value = 3666
$ws.Range("N3").Value = -3898
# Row 4
$ws.Range("H4").Value = 672.25
$ws.Range("I4").Value = 849.5
$ws.Range("J4").Value = 495
$ws.Range("K4").Value = 849.5
$ws.Range("L4").Value = 495
$ws.Range("M4").Value = -737.5
$ws.Range("N4").Value = -719
# Row 6
$ws.Range("H6").Value = 17468.166
$ws.Range("I6").Value = 1500
$ws.Range("J6").Value = 20661.8
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 20661.8
$ws.Range("M6").Value = -1387
$ws.Range("N6").Value = -20887.8
# Row 9
$ws.Range("H9").Value = 1095
$ws.Range("I9").Value = 1730
$ws.Range("J9").Value = 460
$ws.Range("K9").Value = 1730
$ws.Range("L9").Value = 460
$ws.Range("M9").Value = -1560
$ws.Range("N9").Value = -800
# Row 10
$ws.Range("H10").Value = 337600
$ws.Range("I10").Value = 506000
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 506000
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = -505831
$ws.Range("N10").Value = -1138
# Row 11
$ws.Range("H11").Value = 900
$ws.Range("I11").Value = 750
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 750
$ws.Range("L11").Value = 1200
$ws.Range("M11").Value = -611
$ws.Range("N11").Value = -1478
# Row 12
$ws.Range("H12").Value = 50
$ws.Range("J12").Value = 50
$ws.Range("L12").Value = 50
$ws.Range("N12").Value = -330
# Row 15
$ws.Range("H15").Value = 24666.666
# Row 16
$ws.Range("H16").Value = 17468.166
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 20661.8
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 20661.8
$ws.Range("M16").Value = -1250
$ws.Range("N16").Value = -21161.8
# Row 81
$ws.Range("H81").Value = 24666.666
# Row 84
$ws.Range("H84").Value = 24666.666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2334.08
$ws.Range("J22").Value = 2909.111
$ws.Range("L22").Value = 2909.111
$ws.Range("N22").Value = -3499.111
# Row 27
$ws.Range("H27").Value = 2334.08
$ws.Range("J27").Value = 2909.111
$ws.Range("L27").Value = 2909.111
$ws.Range("N27").Value = -3123.111
# Row 35
$ws.Range("H35").Value = 515
$ws.Range("I35").Value = 515
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 515
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -179
$ws.Range("N35").Value = ""
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 29999.5
$ws.Range("J104").Value = 29999.5
$ws.Range("L104").Value = 29999.5
$ws.Range("N104").Value = -36987.5
